$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of resale numbers for 2024-01-01 21:58:05 ("Monday", week "00")
# Columns A and D look like a date / a number to Excel's auto-detection, so
# force them to be treated as plain text (matching the original rows) and
# then clear the formatting delta so no new cell style gets introduced.
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "2024-01-01"
$ws.Cells.Item(4, 1).ClearFormats()

$ws.Cells.Item(4, 2).Value = "21:58:05"
$ws.Cells.Item(4, 3).Value = "Monday"

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "00"
$ws.Cells.Item(4, 4).ClearFormats()

$ws.Cells.Item(4, 5).Value = 140128
$ws.Cells.Item(4, 6).Value = 142982
$ws.Cells.Item(4, 7).Value = 172180
$ws.Cells.Item(4, 8).Value = 145313
$ws.Cells.Item(4, 9).Value = -1
$ws.Cells.Item(4, 10).Value = 116931
$ws.Cells.Item(4, 11).Value = 223918
$ws.Cells.Item(4, 12).Value = 248122
$ws.Cells.Item(4, 13).Value = 183533
$ws.Cells.Item(4, 14).Value = 109881
$ws.Cells.Item(4, 15).Value = 39619
$ws.Cells.Item(4, 16).Value = 30594
$ws.Cells.Item(4, 17).Value = 71614
$ws.Cells.Item(4, 18).Value = -1
$ws.Cells.Item(4, 19).Value = 42000
$ws.Cells.Item(4, 20).Value = -1
